# JackieK-WorkLog.xlsx - "Add files via upload"
#
# Appends one new work-log entry (row 35) to Sheet1:
#   Date: 2025-03-14 (serial 45730), Hours: 4,
#   Description: "Changed web app security scanning tool : nikto, tested
#   Ollama tool for analysis report."
#
# The last data row of the table (row 34) previously carried the
# "end of table" cell formatting (style applied to B34/C34). When a new
# row is appended below it, that special formatting moves down to the new
# last row (35) and row 34 reverts to the regular interior-row formatting
# used by the rest of the table (matching row 33, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Write the new row's data --------------------------------------
$ws.Range("A35").Value = 45730
$ws.Range("B35").Value = 4
$ws.Range("C35").Value = "Changed web app security scanning tool : nikto, tested Ollama tool for analysis report."

# --- 2. Carry the old last-row formatting down onto the new last row ---
$ws.Range("A34:C34").Copy()
$ws.Range("A35:C35").PasteSpecial(-4122)   # xlPasteFormats

# --- 3. Restore row 34 to the standard interior-row formatting ---------
$ws.Range("A33:C33").Copy()
$ws.Range("A34:C34").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- 4. Update the sheet's active selection -----------------------------
$ws.Range("C31").Select()
